$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '66.098.54'
$c.Style = $origStyle
$ws.Range("E2").Value = '  -0.25%  '
$c = $ws.Range("D3")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.257.38'
$c.Style = $origStyle
$ws.Range("E3").Value = '  +2.22%  '
$c = $ws.Range("D4")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = $origStyle
$ws.Range("E4").Value = '  +0.01%  '
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '608.13'
$c.Style = $origStyle
$ws.Range("E5").Value = '  +0.44%  '
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '157.27'
$c.Style = $origStyle
$ws.Range("E6").Value = '  +1.65%  '
$ws.Range("E7").Value = '  +0.06%  '
$c = $ws.Range("D8")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.258.54'
$c.Style = $origStyle
$ws.Range("E8").Value = '  +2.32%  '
$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.540'
$c.Style = $origStyle
$ws.Range("E9").Value = '  -0.82%  '
$c = $ws.Range("D10")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.160'
$c.Style = $origStyle
$ws.Range("E10").Value = '  +0.54%  '
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.70'
$c.Style = $origStyle
$ws.Range("E11").Value = '  +0.71%  '
$c = $ws.Range("D12")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.494'
$c.Style = $origStyle
$ws.Range("E12").Value = '  -2.80%  '
$c = $ws.Range("D13")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0000267'
$c.Style = $origStyle
$ws.Range("E13").Value = '  +0.29%  '
$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '38.55'
$c.Style = $origStyle
$ws.Range("E14").Value = '  +0.32%  '
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.779.58'
$c.Style = $origStyle
$ws.Range("E15").Value = '  +1.97%  '
$c = $ws.Range("D16")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '66.075.56'
$c.Style = $origStyle
$ws.Range("E16").Value = '  -0.22%  '
$c = $ws.Range("D17")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.248.54'
$c.Style = $origStyle
$ws.Range("E17").Value = '  +1.98%  '
$c = $ws.Range("D18")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.30'
$c.Style = $origStyle
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("E19").Value = '  +1.21%  '
$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '498.34'
$c.Style = $origStyle
$ws.Range("E20").Value = '  -2.15%  '
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '15.30'
$c.Style = $origStyle
$ws.Range("E21").Value = '  +0.06%  '
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.747'
$c.Style = $origStyle
$ws.Range("E22").Value = '  +2.14%  '
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '8.04'
$c.Style = $origStyle
$ws.Range("E23").Value = '  -0.06%  '
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '14.57'
$c.Style = $origStyle
$ws.Range("E24").Value = '  -1.73%  '
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '86.89'
$c.Style = $origStyle
$ws.Range("E25").Value = '  +2.81%  '
$ws.Range("E26").Value = '  -0.09%  '
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.02'
$c.Style = $origStyle
$ws.Range("E27").Value = '  +0.87%  '
$c = $ws.Range("D28")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.08'
$c.Style = $origStyle
$ws.Range("E28").Value = '  -1.08%  '
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.36'
$c.Style = $origStyle
$ws.Range("E29").Value = '  -1.58%  '
$c = $ws.Range("D30")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.132'
$c.Style = $origStyle
$ws.Range("E30").Value = '  +46.69%  '
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.07'
$c.Style = $origStyle
$ws.Range("E31").Value = '  +0.51%  '
$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.84'
$c.Style = $origStyle
$ws.Range("E32").Value = '  -7.09%  '
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '27.86'
$c.Style = $origStyle
$ws.Range("E33").Value = '  -0.44%  '
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = $origStyle
$ws.Range("E34").Value = '  -0.13%  '
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.14'
$c.Style = $origStyle
$ws.Range("E35").Value = '  -3.82%  '
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.38'
$c.Style = $origStyle
$ws.Range("E36").Value = '  -1.65%  '
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.40'
$c.Style = $origStyle
$ws.Range("E37").Value = '  +16.51%  '
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '55.71'
$c.Style = $origStyle
$ws.Range("E38").Value = '  +0.36%  '
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '492.98'
$c.Style = $origStyle
$ws.Range("E39").Value = '  -3.58%  '
$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0₃0770'
$c.Style = $origStyle
$ws.Range("E40").Value = '  +6.30%  '
$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0420'
$c.Style = $origStyle
$ws.Range("E41").Value = '  +0.56%  '
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.130'
$c.Style = $origStyle
$ws.Range("E42").Value = '  +1.87%  '
$c = $ws.Range("D43")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '8.76'
$c.Style = $origStyle
$ws.Range("E43").Value = '  -0.26%  '
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.53'
$c.Style = $origStyle
$ws.Range("E44").Value = '  +2.09%  '
$c = $ws.Range("D45")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.994.89'
$c.Style = $origStyle
$ws.Range("E45").Value = '  +5.71%  '
$c = $ws.Range("D46")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.291'
$c.Style = $origStyle
$ws.Range("E46").Value = '  -2.68%  '
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '28.61'
$c.Style = $origStyle
$ws.Range("E47").Value = '  +2.07%  '
$c = $ws.Range("D48")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.46'
$c.Style = $origStyle
$ws.Range("E48").Value = '  +3.16%  '
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.119'
$c.Style = $origStyle
$ws.Range("E49").Value = '  +2.13%  '
$ws.Range("E50").Value = '  +0.01%  '
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '121.26'
$c.Style = $origStyle
$ws.Range("E51").Value = '  -1.70%  '